$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Betfair back/lay odds for rows 3-7 (2026-01-05 fixtures)
$updates = @{
    3 = @{ "K"=4.1; "L"=1.01; "M"=1.05; "N"=1.99; "O"=1.28; "Q"=1.82; "R"=1.16; "S"=1.82; "T"=1.01; "U"=1.01; "V"=1.24; "W"=1.92; "X"=24; "Y"=24; "Z"=48; "AA"=1000; "AB"=14.5; "AC"=12.5; "AD"=25; "AE"=75; "AF"=18.5; "AG"=15.5; "AH"=26; "AI"=85; "AJ"=34; "AK"=30; "AL"=50; "AM"=1000; "AN"=1000; "AO"=1000 }
    4 = @{ "L"=1.01; "M"=1.01; "N"=6.4; "O"=1.14; "Q"=1.44; "R"=1.75; "S"=1.91; "T"=1.66; "U"=2.18; "V"=1.11; "W"=3.05; "X"=42; "Y"=46; "Z"=90; "AA"=250; "AB"=16.5; "AC"=16.5; "AD"=36; "AE"=110; "AF"=14; "AG"=13.5; "AH"=26; "AI"=90; "AJ"=16.5; "AK"=17; "AL"=34; "AM"=100; "AN"=5.3; "AO"=95 }
    5 = @{ "F"=2.1; "K"=3.7; "L"=1.01; "M"=1.01; "N"=1.7; "O"=1.27; "R"=1.08; "S"=2.1; "T"=1.01; "U"=1.01; "V"=1.34; "W"=1.74; "X"=1000; "Y"=1000; "Z"=1000; "AA"=1000; "AB"=1000; "AC"=1000; "AD"=1000; "AE"=1000; "AF"=1000; "AG"=1000; "AH"=1000; "AI"=1000; "AJ"=1000; "AK"=1000; "AL"=1000; "AM"=1000; "AN"=1000; "AO"=1000 }
    6 = @{ "G"=1.93; "H"=5; "I"=6; "J"=3.5; "K"=3.8; "L"=1.01; "M"=1.01; "N"=1.61; "O"=1.01; "P"=1.61; "Q"=2.18; "R"=1.08; "S"=2.18; "T"=1.01; "U"=1.01; "V"=1.2; "W"=2.06; "X"=1000; "Y"=1000; "Z"=1000; "AA"=1000; "AB"=1000; "AC"=1000; "AD"=1000; "AE"=1000; "AF"=1000; "AG"=1000; "AH"=1000; "AI"=1000; "AJ"=1000; "AK"=1000; "AL"=1000; "AM"=1000; "AN"=1000; "AO"=1000 }
    7 = @{ "I"=7.4; "J"=3.65; "K"=4.3; "L"=1.01; "M"=1.01; "N"=1.94; "O"=1.01; "P"=1.9; "Q"=1.76; "R"=1.12; "S"=1.76; "T"=1.01; "U"=1.01; "V"=1.19; "W"=2.2; "X"=1000; "Y"=1000; "Z"=1000; "AA"=1000; "AB"=1000; "AC"=1000; "AD"=1000; "AE"=1000; "AF"=1000; "AG"=1000; "AH"=1000; "AI"=1000; "AJ"=1000; "AK"=1000; "AL"=1000; "AM"=1000; "AN"=1000; "AO"=1000 }
}

foreach ($row in $updates.Keys) {
    foreach ($col in $updates[$row].Keys) {
        $ws.Range("$col$row").Value = $updates[$row][$col]
    }
}
